# Update "want to go" counts (column F) in the "展览" sheet and the
# mirrored "全部类型" sheet (rows are offset by one between the two).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 (sheet1) updates
$wsExhibit.Range("F6").Value  = 10075
$wsExhibit.Range("F8").Value  = 910
$wsExhibit.Range("F10").Value = 6057
$wsExhibit.Range("F12").Value = 296
$wsExhibit.Range("F15").Value = 3086
$wsExhibit.Range("F17").Value = 299

# 全部类型 (sheet4) updates
$wsAll.Range("F7").Value  = 10075
$wsAll.Range("F9").Value  = 910
$wsAll.Range("F11").Value = 6057
$wsAll.Range("F13").Value = 297
$wsAll.Range("F16").Value = 3086
$wsAll.Range("F18").Value = 299

"ok"
